$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5888.222
$ws.Range("I18").Value = 8079.4
$ws.Range("J18").Value = 3149.25
$ws.Range("K18").Value = 8079.4
$ws.Range("L18").Value = 3149.25
$ws.Range("M18").Value = -7795.4
$ws.Range("N18").Value = -3717.25
$ws.Range("H80").Value = 32353.438
$ws.Range("I80").Value = 17184.5
$ws.Range("J80").Value = 41454.8
$ws.Range("K80").Value = 51553.5
$ws.Range("L80").Value = 124364.4
$ws.Range("M80").Value = -50555.5
$ws.Range("N80").Value = -126360.4
$ws.Range("H83").Value = 32353.438
$ws.Range("I83").Value = 17184.5
$ws.Range("J83").Value = 41454.8
$ws.Range("K83").Value = 154660.5
$ws.Range("L83").Value = 373093.2
$ws.Range("M83").Value = -149668.5
$ws.Range("N83").Value = -383077.2
$ws.Range("H92").Value = 886.2308
$ws.Range("I92").Value = 511.0909
$ws.Range("J92").Value = 2949.5
$ws.Range("K92").Value = 511.0909
$ws.Range("L92").Value = 2949.5
$ws.Range("M92").Value = 736.9091000000001
$ws.Range("N92").Value = -5445.5
$ws.Range("H98").Value = 52636020
$ws.Range("I98").Value = 55559640
$ws.Range("J98").Value = 10900
$ws.Range("K98").Value = 55559640
$ws.Range("L98").Value = 10900
$ws.Range("M98").Value = -55558142
$ws.Range("N98").Value = -13896
$ws.Range("H118").Value = 4536.4443
$ws.Range("I118").Value = 4247
$ws.Range("J118").Value = 5549.5
$ws.Range("K118").Value = 12741
$ws.Range("L118").Value = 16648.5
$ws.Range("M118").Value = -11084
$ws.Range("N118").Value = -19962.5
$ws.Range("H122").Value = 52636020
$ws.Range("I122").Value = 55559640
$ws.Range("J122").Value = 10900
$ws.Range("K122").Value = 166678920
$ws.Range("L122").Value = 32700
$ws.Range("M122").Value = -166676470
$ws.Range("N122").Value = -37600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 10432367
$ws.Range("I97").Value = 1346.5
$ws.Range("J97").Value = 20863388
$ws.Range("K97").Value = 1346.5
$ws.Range("L97").Value = 20863388
$ws.Range("M97").Value = -850.5
$ws.Range("N97").Value = -20864380
$ws.Range("H122").Value = 11007.259
$ws.Range("I122").Value = 11383.208
$ws.Range("J122").Value = 7999.6665
$ws.Range("K122").Value = 34149.624
$ws.Range("L122").Value = 23998.9995
$ws.Range("M122").Value = -31699.624
$ws.Range("N122").Value = -28898.9995
$ws.Range("H132").Value = 6556.5454
$ws.Range("I132").Value = 5071.2
$ws.Range("J132").Value = 9739.429
$ws.Range("K132").Value = 15213.6
$ws.Range("L132").Value = 29218.287
$ws.Range("M132").Value = -12683.6
$ws.Range("N132").Value = -34278.287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 59375.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 59375.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 59375.5
$ws.Range("N51").Value = -60357.5
$ws.Range("H86").Value = 47622788
$ws.Range("I86").Value = 2973.2942
$ws.Range("J86").Value = 250006990
$ws.Range("K86").Value = 2973.2942
$ws.Range("L86").Value = 250006990
$ws.Range("M86").Value = -1850.2942
$ws.Range("N86").Value = -250009236
$ws.Range("H89").Value = 47622788
$ws.Range("I89").Value = 2973.2942
$ws.Range("J89").Value = 250006990
$ws.Range("K89").Value = 14866.471
$ws.Range("L89").Value = 1250034950
$ws.Range("M89").Value = -9250.471
$ws.Range("N89").Value = -1250046182
$ws.Range("H94").Value = 965.5806
$ws.Range("I94").Value = 382.86365
$ws.Range("J94").Value = 2390
$ws.Range("K94").Value = 382.86365
$ws.Range("L94").Value = 2390
$ws.Range("M94").Value = 68.13634999999999
$ws.Range("N94").Value = -3292
$ws.Range("H134").Value = 4847.125
$ws.Range("I134").Value = 2927.5417
$ws.Range("J134").Value = 10605.875
$ws.Range("K134").Value = 8782.625100000001
$ws.Range("L134").Value = 31817.625
$ws.Range("M134").Value = -6247.625100000001
$ws.Range("N134").Value = -36887.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 500
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -726
$ws.Range("H22").Value = 433.33334
$ws.Range("I22").Value = 433.33334
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 433.33334
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -83.33334000000002
$ws.Range("H31").Value = 6971.0625
$ws.Range("I31").Value = 2952.0588
$ws.Range("J31").Value = 11525.934
$ws.Range("K31").Value = 2952.0588
$ws.Range("L31").Value = 11525.934
$ws.Range("M31").Value = -2657.0588
$ws.Range("N31").Value = -12115.934
$ws.Range("H34").Value = 6971.0625
$ws.Range("I34").Value = 2952.0588
$ws.Range("J34").Value = 11525.934
$ws.Range("K34").Value = 2952.0588
$ws.Range("L34").Value = 11525.934
$ws.Range("M34").Value = -2750.0588
$ws.Range("N34").Value = -11929.934
$ws.Range("H58").Value = 13163928
$ws.Range("I58").Value = 26317882
$ws.Range("J58").Value = 9974.368
$ws.Range("K58").Value = 26317882
$ws.Range("L58").Value = 9974.368
$ws.Range("M58").Value = -26317679
$ws.Range("N58").Value = -10380.368
$ws.Range("H62").Value = 7258
$ws.Range("I62").Value = 9000
$ws.Range("J62").Value = 6967.6665
$ws.Range("K62").Value = 9000
$ws.Range("L62").Value = 6967.6665
$ws.Range("M62").Value = -8376
$ws.Range("N62").Value = -8215.666499999999
$ws.Range("H65").Value = 7258
$ws.Range("I65").Value = 9000
$ws.Range("J65").Value = 6967.6665
$ws.Range("K65").Value = 45000
$ws.Range("L65").Value = 34838.3325
$ws.Range("M65").Value = -41880
$ws.Range("N65").Value = -41078.3325
$ws.Range("H134").Value = 22230592
$ws.Range("I134").Value = 8341.817999999999
$ws.Range("J134").Value = 43486660
$ws.Range("K134").Value = 25025.454
$ws.Range("L134").Value = 130459980
$ws.Range("M134").Value = -22490.454
$ws.Range("N134").Value = -130465050
$ws.Range("H136").Value = 13163928
$ws.Range("I136").Value = 26317882
$ws.Range("J136").Value = 9974.368
$ws.Range("K136").Value = 78953646
$ws.Range("L136").Value = 29923.104
$ws.Range("M136").Value = -78951096
$ws.Range("N136").Value = -35023.104

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 5333.3335
$ws.Range("I119").Value = 5333.3335
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 16000.0005
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -11162.0005
$ws.Range("H132").Value = 19398.8
$ws.Range("I132").Value = 13999.5
$ws.Range("J132").Value = 22998.334
$ws.Range("K132").Value = 125995.5
$ws.Range("L132").Value = 206985.006
$ws.Range("M132").Value = -123465.5
$ws.Range("N132").Value = -212045.006
$ws.Range("H136").Value = 2030
$ws.Range("I136").Value = 2030
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6090
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -990
$ws.Range("H138").Value = 80653
$ws.Range("I138").Value = 94544.45
$ws.Range("J138").Value = 4250
$ws.Range("K138").Value = 283633.35
$ws.Range("L138").Value = 12750
$ws.Range("M138").Value = -278493.35
$ws.Range("N138").Value = -23030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3858.5144
$ws.Range("I102").Value = 3551.3076
$ws.Range("J102").Value = 4746
$ws.Range("K102").Value = 3551.3076
$ws.Range("L102").Value = 4746
$ws.Range("M102").Value = -1929.3076
$ws.Range("N102").Value = -7990
$ws.Range("H113").Value = 6270.7095
$ws.Range("I113").Value = 3693.25
$ws.Range("J113").Value = 7167.2173
$ws.Range("K113").Value = 3693.25
$ws.Range("L113").Value = 7167.2173
$ws.Range("M113").Value = -1523.25
$ws.Range("N113").Value = -11507.2173
$ws.Range("H122").Value = 64444.47
$ws.Range("I122").Value = 253099.25
$ws.Range("J122").Value = 6396.846
$ws.Range("K122").Value = 759297.75
$ws.Range("L122").Value = 19190.538
$ws.Range("M122").Value = -756847.75
$ws.Range("N122").Value = -24090.538
$ws.Range("H126").Value = 2903.0625
$ws.Range("I126").Value = 2909.5
$ws.Range("J126").Value = 2899.2
$ws.Range("K126").Value = 8728.5
$ws.Range("L126").Value = 8697.599999999999
$ws.Range("M126").Value = -6258.5
$ws.Range("N126").Value = -13637.6
$ws.Range("H134").Value = 99998
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 99998
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 299994
$ws.Range("N134").Value = -305064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7361.385
$ws.Range("I7").Value = 5896.5
$ws.Range("J7").Value = 7627.727
$ws.Range("K7").Value = 5896.5
$ws.Range("L7").Value = 7627.727
$ws.Range("M7").Value = -5784.5
$ws.Range("N7").Value = -7851.727
$ws.Range("H80").Value = 49995
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 49995
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 49995
$ws.Range("N80").Value = -52241
$ws.Range("H82").Value = 2203.25
$ws.Range("I82").Value = 2026
$ws.Range("J82").Value = 2451.4
$ws.Range("K82").Value = 2026
$ws.Range("L82").Value = 2451.4
$ws.Range("M82").Value = -1665
$ws.Range("N82").Value = -3173.4
$ws.Range("H83").Value = 49995
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 49995
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 149985
$ws.Range("N83").Value = -161217
$ws.Range("H85").Value = 2203.25
$ws.Range("I85").Value = 2026
$ws.Range("J85").Value = 2451.4
$ws.Range("K85").Value = 2026
$ws.Range("L85").Value = 2451.4
$ws.Range("M85").Value = -778
$ws.Range("N85").Value = -4947.4
$ws.Range("H126").Value = 7361.385
$ws.Range("I126").Value = 5896.5
$ws.Range("J126").Value = 7627.727
$ws.Range("K126").Value = 17689.5
$ws.Range("L126").Value = 22883.181
$ws.Range("M126").Value = -15219.5
$ws.Range("N126").Value = -27823.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 24998
$ws.Range("I15").Value = 24998
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 24998
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -24710
$ws.Range("H96").Value = 2059.6667
$ws.Range("I96").Value = 2059.6667
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2059.6667
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -686.6667000000002
$ws.Range("N96").ClearContents()
$ws.Range("H122").Value = 4326.1904
$ws.Range("I122").Value = 2245.75
$ws.Range("J122").Value = 7100.1113
$ws.Range("K122").Value = 6737.25
$ws.Range("L122").Value = 21300.3339
$ws.Range("M122").Value = -4287.25
$ws.Range("N122").Value = -26200.3339
$ws.Range("H126").Value = 3004.5557
$ws.Range("I126").Value = 1990.7693
$ws.Range("J126").Value = 3945.9285
$ws.Range("K126").Value = 5972.3079
$ws.Range("L126").Value = 11837.7855
$ws.Range("M126").Value = -3502.3079
$ws.Range("N126").Value = -16777.7855
